$wb = $excel.ActiveWorkbook

# Sheet "NBR" - update Reaction_number column (C2:C20)
$wsNBR = $wb.Worksheets.Item("NBR")
$nbrValues = @(700,696,690,683,679,669,663,662,655,654,645,598,597,588,574,568,559,547,544)
for ($i = 0; $i -lt $nbrValues.Length; $i++) {
    $row = $i + 2
    $wsNBR.Cells.Item($row, 3).Value = $nbrValues[$i]
}

# Sheet "BAR" - update Reaction_number column (C2:C20)
$wsBAR = $wb.Worksheets.Item("BAR")
$barValues = @(629,627,626,630,629,630,639,627,629,631,621,621,622,625,626,621,620,619,619)
for ($i = 0; $i -lt $barValues.Length; $i++) {
    $row = $i + 2
    $wsBAR.Cells.Item($row, 3).Value = $barValues[$i]
}
